$d = $word.ActiveDocument

# Locate the target paragraph (starts with "Ovaj predlo..." - the long
# Croatian explanatory paragraph that gets collapsed down to the
# "${description}" placeholder).
$target = $null
for ($i = 1; $i -le $d.Paragraphs.Count; $i++) {
    $cand = $d.Paragraphs($i)
    if ($cand.Range.Text.StartsWith("Ovaj predlo")) {
        $target = $cand
        break
    }
}

$pStart = $target.Range.Start
$pEnd = $target.Range.End

# Address the paragraph's whole body (excluding the trailing paragraph
# mark) and replace it in one shot with the "${description}" placeholder,
# split into three runs with proofErr spell-check markers bracketing the
# word "description" - mirroring the sibling "${block}" placeholder
# paragraph elsewhere in this document.
$bodyRange = $d.Range($pStart, $pEnd - 1)

$xmlFrag = '<?xml version="1.0" encoding="UTF-8" standalone="yes"?><pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage"><pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml"><pkg:xmlData><w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:body><w:p><w:r><w:rPr><w:rFonts w:cs="Times New Roman"/><w:szCs w:val="24"/><w:lang w:val="hr-HR"/></w:rPr><w:t>${</w:t></w:r><w:proofErr w:type="spellStart"/><w:r><w:rPr><w:rFonts w:cs="Times New Roman"/><w:szCs w:val="24"/><w:lang w:val="hr-HR"/></w:rPr><w:t>description</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r><w:rPr><w:rFonts w:cs="Times New Roman"/><w:szCs w:val="24"/><w:lang w:val="hr-HR"/></w:rPr><w:t>}</w:t></w:r></w:p></w:body></w:document></pkg:xmlData></pkg:part></pkg:package>'
$bodyRange.InsertXML($xmlFrag)

$d.Save()
